$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 350, pushing the existing rows 350-364 down to 351-365.
$ws.Rows.Item(350).Insert()

# Populate the newly inserted row 350 with the new record.
$ws.Cells.Item(350, 1).Value = 4
$ws.Cells.Item(350, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(350, 3).Value = "Los Lagos"
$ws.Cells.Item(350, 4).Value = 45041
$ws.Cells.Item(350, 5).Value = 10
$ws.Cells.Item(350, 6).Value = "Fruta"
$ws.Cells.Item(350, 7).Value = 100108
$ws.Cells.Item(350, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(350, 9).Value = 100108002
$ws.Cells.Item(350, 10).Value = "Mango"
$ws.Cells.Item(350, 11).Value = "Sin especificar"
$ws.Cells.Item(350, 12).Value = "Primera"
$ws.Cells.Item(350, 13).Value = 200
$ws.Cells.Item(350, 14).Value = 8000
$ws.Cells.Item(350, 15).Value = 8500
$ws.Cells.Item(350, 16).Value = 8250
$ws.Cells.Item(350, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(350, 18).Value = "Perú"
$ws.Cells.Item(350, 19).Value = 2062
$ws.Cells.Item(350, 20).Value = 4
